# Generate Report for Handoff
# Replaces the prior handoff's file id/hash ("a4b82b49-...") with the new
# handoff's id/hash ("5fefe0d4-...") across the Overview / zh-cn / de-de
# sheets, refreshes the handoff timestamps, and clears the (not-yet-produced)
# "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns on the locale sheets.

$wb = $excel.ActiveWorkbook

$oldId = "a4b82b49-f1f5-468d-951b-ab75af631667"
$newId = "5fefe0d4-645d-458a-ac26-ee5a24112332"

$oldZhXlf = "$oldId.33860ba8260f646f8e82229ac1933ae4acd1a790.zh-cn.xlf"
$newZhXlf = "$newId.a03fef4c6a54f19a1527384072b7a05b202b28ae.zh-cn.xlf"
$oldDeXlf = "$oldId.33860ba8260f646f8e82229ac1933ae4acd1a790.de-de.xlf"
$newDeXlf = "$newId.a03fef4c6a54f19a1527384072b7a05b202b28ae.de-de.xlf"

$mdName = "$newId.md"
$mdDisplayOverview = "e2e\$newId.md"

$ghAddr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/53715ad08294096f446025a82d1c6680dd843c51/e2e/$mdName"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $mdName
$wsOverview.Range("B2").Value = $mdDisplayOverview
$wsOverview.Range("G2").Value = "2016-09-04 09:03:48"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $ghAddr, "", "", $mdDisplayOverview)

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $mdName
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = "2016-09-04 09:03:43"
$wsZh.Range("I2").Value = ""
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = "0001-01-01 00:00:00"

$wsZh.Range("I2").Style = "Normal"

$wsZh.Columns.Item(9).ColumnWidth = 17.8
$wsZh.Columns.Item(10).ColumnWidth = 20.8

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $ghAddr, "", "", $mdName)

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $mdName
$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = "2016-09-04 09:03:48"
$wsDe.Range("I2").Value = ""
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDe.Range("I2").Style = "Normal"

$wsDe.Columns.Item(9).ColumnWidth = 17.8
$wsDe.Columns.Item(10).ColumnWidth = 20.8

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $ghAddr, "", "", $mdName)
